# This workbook contains a weekly price feed for "Pimiento" (bell pepper) at
# "Terminal Hortofrutícola Agro Chillán". A new week of data is published at
# the top of this product's block (rows 491-538), which pushes the existing
# rows down by two positions (491-538 -> 493-540).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 491, shifting rows 491:538 down to 493:540.
$ws.Rows.Item(491).Insert()
$ws.Rows.Item(491).Insert()

# New row 491: Zafiro rojo, Primera
$ws.Range("A491").Value = 7
$ws.Range("B491").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C491").Value = "Ñuble"
$ws.Range("D491").Value = 45132
$ws.Range("E491").Value = 16
$ws.Range("F491").Value = 100112002
$ws.Range("G491").Value = "Pimiento"
$ws.Range("H491").Value = "Zafiro rojo"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 80
$ws.Range("K491").Value = 14000
$ws.Range("L491").Value = 14000
$ws.Range("M491").Value = 14000
$ws.Range("N491").Value = "$/caja 15 kilos"
$ws.Range("O491").Value = "Región de Arica y Parinacota"
$ws.Range("P491").Value = 933
$ws.Range("Q491").Value = 15
$ws.Range("R491").Value = "Hortaliza"

# New row 492: Zafiro verde, Primera
$ws.Range("A492").Value = 7
$ws.Range("B492").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C492").Value = "Ñuble"
$ws.Range("D492").Value = 45132
$ws.Range("E492").Value = 16
$ws.Range("F492").Value = 100112002
$ws.Range("G492").Value = "Pimiento"
$ws.Range("H492").Value = "Zafiro verde"
$ws.Range("I492").Value = "Primera"
$ws.Range("J492").Value = 80
$ws.Range("K492").Value = 13000
$ws.Range("L492").Value = 13000
$ws.Range("M492").Value = 13000
$ws.Range("N492").Value = "$/caja 15 kilos"
$ws.Range("O492").Value = "Región de Arica y Parinacota"
$ws.Range("P492").Value = 867
$ws.Range("Q492").Value = 15
$ws.Range("R492").Value = "Hortaliza"
